# Scheduled-runner price/profit refresh across the Leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice /
# NQ / HQ price & profit columns (H,I,J,K,L,M,N) for the rows whose
# underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1119.6
$ws.Range("I6").Value = 1119.6
$ws.Range("K6").Value = 3358.8
$ws.Range("M6").Value = -3246.8

$ws.Range("H8").Value = 28044.889
$ws.Range("I8").Value = 166775.17
$ws.Range("K8").Value = 500325.51
$ws.Range("M8").Value = -500186.51

$ws.Range("H17").Value = 2077.7334
$ws.Range("J17").Value = 2108
$ws.Range("L17").Value = 6324
$ws.Range("N17").Value = -6660

$ws.Range("H28").Value = 45078.305
$ws.Range("I28").Value = 84496.414
$ws.Range("J28").Value = 2076.7273
$ws.Range("K28").Value = 84496.414
$ws.Range("L28").Value = 2076.7273
$ws.Range("M28").Value = -84011.414
$ws.Range("N28").Value = -3046.7273

$ws.Range("H33").Value = 147.3077
$ws.Range("I33").Value = 142.75
$ws.Range("K33").Value = 142.75
$ws.Range("M33").Value = 86.25

$ws.Range("H39").Value = 171.5238
$ws.Range("J39").Value = 232.55556
$ws.Range("L39").Value = 697.66668
$ws.Range("N39").Value = -1289.66668

$ws.Range("H64").Value = 7514.2856
$ws.Range("I64").Value = 3533.3333
$ws.Range("J64").Value = 10500
$ws.Range("K64").Value = 3533.3333
$ws.Range("L64").Value = 10500
$ws.Range("M64").Value = -3285.3333
$ws.Range("N64").Value = -10996

$ws.Range("H67").Value = 7514.2856
$ws.Range("I67").Value = 3533.3333
$ws.Range("J67").Value = 10500
$ws.Range("K67").Value = 3533.3333
$ws.Range("L67").Value = 10500
$ws.Range("M67").Value = -2675.3333
$ws.Range("N67").Value = -12216

$ws.Range("H132").Value = 1803.3877
$ws.Range("I132").Value = 1034.15
$ws.Range("J132").Value = 5222.222
$ws.Range("K132").Value = 3102.45
$ws.Range("L132").Value = 15666.666
$ws.Range("M132").Value = -572.4500000000003
$ws.Range("N132").Value = -20726.666

$ws.Range("H137").Value = 2420.8508
$ws.Range("I137").Value = 1213.9487
$ws.Range("K137").Value = 3641.8461
$ws.Range("M137").Value = -1091.8461

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2594.9
$ws.Range("I32").Value = 2205.1973
$ws.Range("K32").Value = 2205.1973
$ws.Range("M32").Value = -1918.1973

$ws.Range("H61").Value = 2600.739
$ws.Range("J61").Value = 6749.375
$ws.Range("L61").Value = 6749.375
$ws.Range("N61").Value = -7173.375

$ws.Range("H74").Value = 1759.1613
$ws.Range("I74").Value = 1415.6072
$ws.Range("K74").Value = 1415.6072
$ws.Range("M74").Value = -541.6071999999999

$ws.Range("H77").Value = 1759.1613
$ws.Range("I77").Value = 1415.6072
$ws.Range("K77").Value = 7078.036
$ws.Range("M77").Value = -2710.036

$ws.Range("H132").Value = 3314.2083
$ws.Range("I132").Value = 1152.7115
$ws.Range("J132").Value = 8934.1
$ws.Range("K132").Value = 3458.1345
$ws.Range("L132").Value = 26802.3
$ws.Range("M132").Value = -928.1344999999997
$ws.Range("N132").Value = -31862.3

$ws.Range("H135").Value = 98214.5
$ws.Range("J135").Value = 98214.5
$ws.Range("L135").Value = 98214.5
$ws.Range("N135").Value = -108354.5

$ws.Range("H136").Value = 2600.739
$ws.Range("J136").Value = 6749.375
$ws.Range("L136").Value = 20248.125
$ws.Range("N136").Value = -25348.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 214592
$ws.Range("J42").Value = 214592
$ws.Range("L42").Value = 214592
$ws.Range("N42").Value = -215248

$ws.Range("H134").Value = 4061.8086
$ws.Range("I134").Value = 2769.0667
$ws.Range("J134").Value = 6343.1177
$ws.Range("K134").Value = 8307.2001
$ws.Range("L134").Value = 19029.3531
$ws.Range("M134").Value = -5772.2001
$ws.Range("N134").Value = -24099.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 211459.25
$ws.Range("I58").Value = 334744.38
$ws.Range("J58").Value = 5984.0557
$ws.Range("K58").Value = 334744.38
$ws.Range("L58").Value = 5984.0557
$ws.Range("M58").Value = -334541.38
$ws.Range("N58").Value = -6390.0557

$ws.Range("H86").Value = 8105.3335
$ws.Range("I86").Value = 7628.4287
$ws.Range("J86").Value = 8773
$ws.Range("K86").Value = 7628.4287
$ws.Range("L86").Value = 8773
$ws.Range("M86").Value = -6505.4287
$ws.Range("N86").Value = -11019

$ws.Range("H89").Value = 8105.3335
$ws.Range("I89").Value = 7628.4287
$ws.Range("J89").Value = 8773
$ws.Range("K89").Value = 38142.14350000001
$ws.Range("L89").Value = 43865
$ws.Range("M89").Value = -32526.14350000001
$ws.Range("N89").Value = -55097

$ws.Range("H94").Value = 803.4545000000001
$ws.Range("J94").Value = 794.35
$ws.Range("L94").Value = 794.35
$ws.Range("N94").Value = -1696.35

$ws.Range("H136").Value = 211459.25
$ws.Range("I136").Value = 334744.38
$ws.Range("J136").Value = 5984.0557
$ws.Range("K136").Value = 1004233.14
$ws.Range("L136").Value = 17952.1671
$ws.Range("M136").Value = -1001683.14
$ws.Range("N136").Value = -23052.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 109.333336
$ws.Range("I10").Value = 39
$ws.Range("K10").Value = 117
$ws.Range("M10").Value = 22

$ws.Range("H13").Value = 1913.3636
$ws.Range("I13").Value = 721.4286
$ws.Range("J13").Value = 3999.25
$ws.Range("K13").Value = 2164.2858
$ws.Range("L13").Value = 11997.75
$ws.Range("M13").Value = -1996.2858
$ws.Range("N13").Value = -12333.75

$ws.Range("H23").Value = 212
$ws.Range("J23").Value = 207.77777
$ws.Range("L23").Value = 623.33331
$ws.Range("N23").Value = -1093.33331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 43998
$ws.Range("J123").Value = 43998
$ws.Range("L123").Value = 43998
$ws.Range("N123").Value = -48898

$ws.Range("H132").Value = 439084.7
$ws.Range("I132").Value = 557997.2
$ws.Range("K132").Value = 1673991.6
$ws.Range("M132").Value = -1671461.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1849.4706
$ws.Range("J22").Value = 1774.4445
$ws.Range("L22").Value = 1774.4445
$ws.Range("N22").Value = -2364.4445

$ws.Range("H27").Value = 1849.4706
$ws.Range("J27").Value = 1774.4445
$ws.Range("L27").Value = 1774.4445
$ws.Range("N27").Value = -1988.4445

$ws.Range("H132").Value = 5522.077
$ws.Range("I132").Value = 4798
$ws.Range("K132").Value = 14394
$ws.Range("M132").Value = -11864

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H62").Value = 7465
$ws.Range("I62").Value = 8749.5
$ws.Range("K62").Value = 8749.5
$ws.Range("M62").Value = -8125.5

$ws.Range("H65").Value = 7465
$ws.Range("I65").Value = 8749.5
$ws.Range("K65").Value = 43747.5
$ws.Range("M65").Value = -40627.5

$ws.Range("H68").Value = 40749.75
$ws.Range("I68").Value = 35000
$ws.Range("K68").Value = 35000
$ws.Range("M68").Value = -34189

$ws.Range("H71").Value = 40749.75
$ws.Range("I71").Value = 35000
$ws.Range("K71").Value = 105000
$ws.Range("M71").Value = -100944

$ws.Range("H113").Value = 1655.4584
$ws.Range("I113").Value = 1197.1578
$ws.Range("J113").Value = 3397
$ws.Range("K113").Value = 3591.4734
$ws.Range("L113").Value = 10191
$ws.Range("M113").Value = -1421.4734
$ws.Range("N113").Value = -14531

$ws.Range("H132").Value = 2533.3416
$ws.Range("I132").Value = 1079.9333
$ws.Range("K132").Value = 3239.7999
$ws.Range("M132").Value = -709.7999

$ws.Range("H136").Value = 2282.4849
$ws.Range("I136").Value = 1656.68
$ws.Range("J136").Value = 4238.125
$ws.Range("K136").Value = 4970.04
$ws.Range("L136").Value = 12714.375
$ws.Range("M136").Value = -2420.04
$ws.Range("N136").Value = -17814.375
